# Apply the "esa to route xlsx and route ids" edit:
#  1. Insert a new column A ("ID") before the existing data, shifting
#     Category..Req reqs from A:I to B:J.
#  2. Fill the new ID column with a running row id (1..30) for each
#     data row.
#  3. Append two new rows describing the new "ESA" category
#     (Add an exercise / Add a stretch).
#  4. Fix up column widths / header style for the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a brand new column before column A; this shifts all existing
#    columns (A..I) one place to the right (B..J) together with their
#    values/styles/column widths.
$ws.Columns("A").Insert()

# 2. Header for the new column.
$ws.Cells.Item(1, 1).Value2 = "ID"

# 3. Running id numbers for the existing 29 data rows (rows 2..29 before
#    the two new rows are appended below).
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $r - 1
}

# 3b. Row 14 ("View one workout from an external user") was missing its
#     Security value; fill it in now that it lines up with the other rows.
$ws.Cells.Item(14, 9).Value2 = "JWT"

# 4. New rows for the "ESA" category.
$ws.Cells.Item(30, 1).Value2 = 29
$ws.Cells.Item(30, 2).Value2 = "ESA"
$ws.Cells.Item(30, 3).Value2 = "Add an exercise"
$ws.Cells.Item(30, 4).Value2 = "POST"
$ws.Cells.Item(30, 4).HorizontalAlignment = -4108
$ws.Cells.Item(30, 5).Value2 = "/exercises"
$ws.Cells.Item(30, 9).Value2 = "spec password"
$ws.Cells.Item(30, 10).Value2 = "[]of datatype.Exercise"

$ws.Cells.Item(31, 1).Value2 = 30
$ws.Cells.Item(31, 2).Value2 = "ESA"
$ws.Cells.Item(31, 3).Value2 = "Add a stretch"
$ws.Cells.Item(31, 4).Value2 = "POST"
$ws.Cells.Item(31, 4).HorizontalAlignment = -4108
$ws.Cells.Item(31, 5).Value2 = "/stretches"
$ws.Cells.Item(31, 9).Value2 = "spec password"
$ws.Cells.Item(31, 10).Value2 = "[]of datatype.Stretch"

# 5. Column widths: the "Security" column (now column I) grew to fit
#    "spec password" and gained a best-fit flag.
$ws.Columns("I").ColumnWidth = 13.140625

# 6. Update the used range / selection to match the new data extent.
$ws.Range("A19").Select() | Out-Null
